$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D
$ws.Range("D1").Value = "GNPXC"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell (A1) onto the new header cell D1.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# New data values for column D (GNPXC), rows 2-26
$values = @(
    -21.7285019360561,
    -25.08060966185799,
    -27.10010129686633,
    -27.06908960397979,
    -28.8825220754527,
    -29.94859058006045,
    -25.89324583683593,
    -19.93608370400428,
    -16.82127840733422,
    -12.86973756267729,
    -8.211208407277148,
    -7.566958904969685,
    -6.956106265842839,
    -4.780545139431,
    -5.055153423506624,
    0.4800652948802238,
    12.27015730682708,
    18.24331659142144,
    21.94851821104215,
    26.64014873022778,
    38.03504935997214,
    45.68189272286769,
    56.21589622837237,
    70.5583284097419,
    84.2374238461634
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}
